$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new year column (T) with its three data values, mirroring the
# style of the existing last column (S).
$ws.Range("T4").Value = 2023
$ws.Range("T5").Value = 1.4
$ws.Range("T6").Value = 8.1999999999999993

# Copy formatting from column S so the new column matches the existing
# table styling (borders, number formats, fonts, etc.)
$ws.Range("S3:S6").Copy() | Out-Null
$ws.Range("T3:T6").PasteSpecial(-4122) | Out-Null

# Make sure the new header cell keeps the integer year format (no decimals)
$ws.Range("T4").Value = 2023

# Adjust column widths: A:C slightly narrower, D:T uniform width matching
# the new data range, and clear the old selection anchor.
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 36.140625
$ws.Range("D1:T1").EntireColumn.ColumnWidth = 8.5703125

# Row heights for header rows 2-4 become an explicit 13.5pt custom height
$ws.Range("A2:A4").EntireRow.RowHeight = 13.5

$ws.Range("A1").Select() | Out-Null
